$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44659
$ws.Range("K2").Value = 2500
$ws.Range("M2").Value = 2722
$ws.Range("P2").Value = 454
$ws.Range("D3").Value = 44637
$ws.Range("J3").Value = 170
$ws.Range("K3").Value = 2800
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2906
$ws.Range("P3").Value = 484
$ws.Range("D4").Value = 44685
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = 3267
$ws.Range("P4").Value = 544
$ws.Range("D5").Value = 44643
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 2800
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 2911
$ws.Range("P5").Value = 485
$ws.Range("D6").Value = 44957
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1857
$ws.Range("P6").Value = 310
$ws.Range("D7").Value = 44630
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 2722
$ws.Range("P7").Value = 454
$ws.Range("D8").Value = 44987
$ws.Range("J8").Value = 130
$ws.Range("K8").Value = 4500
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = 4692
$ws.Range("P8").Value = 782
$ws.Range("D9").Value = 44876
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6812
$ws.Range("P9").Value = 1135
$ws.Range("D10").Value = 44658
$ws.Range("J10").Value = 180
$ws.Range("M10").Value = 2778
$ws.Range("P10").Value = 463
$ws.Range("D11").Value = 45021
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 4500
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 4700
$ws.Range("P11").Value = 783
$ws.Range("D12").Value = 44672
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3500
$ws.Range("M12").Value = 3286
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 548
$ws.Range("D14").Value = 44671
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 3500
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = 3733
$ws.Range("P14").Value = 622
$ws.Range("D15").Value = 44650
$ws.Range("J15").Value = 130
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 3500
$ws.Range("M15").Value = 3308
$ws.Range("P15").Value = 551
$ws.Range("D16").Value = 44644
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 2786
$ws.Range("O16").Value = "Provincia de Chacabuco"
$ws.Range("P16").Value = 464
